$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Quantity Available changed ---
$ws.Range("K3").Value = 233027

# --- Row 4: Quantity Available changed ---
$ws.Range("K4").Value = 701455

# --- Row 8: Quantity + Extended Price changed ---
$ws.Range("H8").Value = 4
$ws.Range("J8").Value = "'`$1.68"
$ws.Range("J8").Style = "Normal"

# --- Row 10: Quantity + Extended Price changed ---
$ws.Range("H10").Value = 1
$ws.Range("J10").Value = "'`$0.92"
$ws.Range("J10").Style = "Normal"

# --- Row 12: new BOM line (UART/USB, JTAG, I2C header) ---
$ws.Range("A12").Value = "'0022232041"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = "Molex, LLC"
$ws.Range("C12").Value = "WM4202-ND"
$ws.Range("D12").Value = "'"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "Bulk"
$ws.Range("G12").Value = "Active"
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 0.28
$ws.Range("J12").Value = "'`$0.56"
$ws.Range("J12").Style = "Normal"
$ws.Range("K12").Value = 51523
$ws.Range("L12").Value = "7 Weeks"
$ws.Range("M12").Value = "CONN HEADER 4POS .100 VERT TIN"
$ws.Range("N12").Value = "RoHS Compliant"
$ws.Range("O12").Value = "Lead free"
$ws.Range("P12").Value = "REACH Unaffected"
